$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlLeft alignment constant (used by the existing "is_active" H column style)
$xlLeft = -4131

# New device rows 157-161 (group "32": Finger Print / IRIS / Web Camera /
# Document Scanner / Printer scanners), following the exact same pattern as
# every previous group of 5 rows in the sheet.
$newRows = @(
    @{ Row=157; Id=3000176; Name="Finger Print Scanner 32"; Mac="80-75-40-E8-CA-24"; Serial="BS563Q2230824"; DspecId=165 },
    @{ Row=158; Id=3000177; Name="IRIS Scanner 32";         Mac="0E-1A-14-4A-6D-3A"; Serial="BS563Q2230825"; DspecId=327 },
    @{ Row=159; Id=3000178; Name="Web Camera 32";           Mac="65-13-7F-0F-F7-53"; Serial="BS563Q2230826"; DspecId=736 },
    @{ Row=160; Id=3000179; Name="Document Scanner 32";     Mac="73-C4-DE-8E-C9-8D"; Serial="BS563Q2230827"; DspecId=801 },
    @{ Row=161; Id=3000180; Name="Printer 32";               Mac="EC-74-AB-E0-0F-38"; Serial="BS563Q2230828"; DspecId=920 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.Id
    $ws.Range("B$row").Value = $r.Name
    $ws.Range("C$row").Value = $r.Mac
    $ws.Range("D$row").Value = $r.Serial
    $ws.Range("F$row").Value = $r.DspecId
    $ws.Range("G$row").Value = "eng"
    $ws.Range("H$row").Value = $true
    $ws.Range("H$row").HorizontalAlignment = $xlLeft
    $ws.Range("I$row").Value = "superadmin"
    $ws.Range("J$row").Value = "now()"
    $ws.Range("K$row").Value = "now()"
}

# Trailing blank rows 162-166 keep the left-aligned style in column H, same
# as the rest of the "is_active" column, but carry no values.
$ws.Range("H162:H166").HorizontalAlignment = $xlLeft

# Update view/selection to match where the user ended up after the edit.
$ws.Range("E159").Select()
